# Split the combined "<exp>...</exp>" abbreviation-expansion markup that
# currently lives inside a single plain run into five runs so the
# "<exp>" / "</exp>" tag text renders in the small grey Courier New
# "markup" style used elsewhere in the document, while the text in
# between (and around) keeps the normal black Arial body style.
#
# There are three occurrences in the document, in this order:
#   1. demeure co<exp>mm</exp>e suspendue sur quelque vuide de lasseurer affin
#   2.  prinse legerem<exp>ent</exp> avecq la poincte dun
#   3. Et encores quelle soit foible co<exp>mm</exp>e la foeuille d

$d = $word.ActiveDocument

# Re-usable: search for $tagText starting at character position
# $searchFrom (to the end of the document), and apply the markup-tag
# character formatting to whatever is found, in place. Returns the
# matched Range so the caller can read its Start/End.
function Format-MarkupTag($searchFrom, $tagText) {
    $rng = $d.Range($searchFrom, $d.Content.End)
    $rng.Find.ClearFormatting()
    $rng.Find.Replacement.ClearFormatting()
    $rng.Find.Replacement.Font.Name = "Courier New"
    $rng.Find.Replacement.Font.NameAscii = "Courier New"
    $rng.Find.Replacement.Font.NameFarEast = "Courier New"
    $rng.Find.Replacement.Font.NameBi = "Courier New"
    $rng.Find.Replacement.Font.Size = 7
    $rng.Find.Replacement.Font.Color = 0xa9a9a9
    $rng.Find.Execute($tagText, $true, $false, $false, $false, $false, $true, 1, $false, $tagText, 2)
    return $rng
}

$searchPos = 0
for ($i = 0; $i -lt 3; $i++) {
    $openTag = Format-MarkupTag $searchPos "<exp>"
    $closeTag = Format-MarkupTag $openTag.End "</exp>"
    $searchPos = $closeTag.End
}
